$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header section (Date / Team Name / Total Number of Team Members) ---
$ws.Range("B3").Value = (Get-Date -Year 2020 -Month 12 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B4").Value = "Limette"
$ws.Range("B5").Value = 4

# --- Team member names & salaries (rows 8-11) ---
$ws.Range("A8").Value = "Lukas Hasler"
$ws.Range("B8").Value = 100
$ws.Range("A9").Value = "Pascal Strebel"
$ws.Range("B9").Value = 100
$ws.Range("A10").Value = "Cedric Weibel"
$ws.Range("B10").Value = 100
$ws.Range("A11").Value = "Robin Schmidiger"
$ws.Range("B11").Value = 100
$ws.Range("A12:B12").ClearContents()

# --- Tasks completed / to complete (rows 19-20) ---
$ws.Range("A19").Value = "Implemented optimization in backend"
$ws.Range("B19").Value = "Finalize communication of front- and backend"
$ws.Range("A20").Value = "Adapted frontend to talk with backend"
$ws.Range("B20").Value = "Prepare final presentation"

# --- Formulas should recalc automatically, but re-assert them to be safe ---
$ws.Range("B14").Formula = "=SUM(B8:B12)"
$ws.Range("B15").Formula = "=B5*100-B14"

# --- Font tweak for the new task rows (10pt instead of default 12pt) ---
$ws.Range("A19:B20").Font.Size = 10

# --- Row 18 height adjusts (auto re-fit by Excel after re-save) ---
$ws.Rows.Item(18).RowHeight = 39

# --- Page setup (paper size / orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection cosmetic change ---
[void]$ws.Range("F15").Select()

Write-Host "done"
